$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Create brand-new cells (format copied from a same-column reference cell,
#     forced to Text so numeric/date-looking strings are kept as literal text) ---
$ws.Range("A13").NumberFormat = "@"
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("A3").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "01/01/2022"
$ws.Range("B2").Copy()
$ws.Range("B15").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "01/01/2022"
$ws.Range("C2").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = "5009972 - Gilberto Carvalho Coelho"
$ws.Range("B2").Copy()
$ws.Range("B18").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "5009972 - Gilberto Carvalho Coelho"
$ws.Range("C2").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("B23").NumberFormat = "@"
$ws.Range("B23").Value = "LOM3013 -  Ciência dos Materiais  (Requisito fraco)`n"
$ws.Range("B2").Copy()
$ws.Range("B23").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = "LOM3013 -  Ciência dos Materiais  (Requisito fraco)`n"
$ws.Range("C2").Copy()
$ws.Range("C23").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Update cells whose text changes (kept Text-typed so values such as dates
#     are not reinterpreted, and style index is preserved) ---
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = "5009972 - Gilberto Carvalho Coelho"
$ws.Range("B2").Copy()
$ws.Range("B10").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "5009972 - Gilberto Carvalho Coelho"
$ws.Range("C2").Copy()
$ws.Range("C10").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "Semestral"
$ws.Range("B2").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "Semestral"
$ws.Range("C2").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A14").NumberFormat = "@"
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("A3").Copy()
$ws.Range("A14").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A15").NumberFormat = "@"
$ws.Range("A15").Value = "Programa:"
$ws.Range("A3").Copy()
$ws.Range("A15").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A16").NumberFormat = "@"
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("A3").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A17").NumberFormat = "@"
$ws.Range("A17").Value = "Avaliação:"
$ws.Range("A3").Copy()
$ws.Range("A17").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A18").NumberFormat = "@"
$ws.Range("A18").Value = "Método:"
$ws.Range("A3").Copy()
$ws.Range("A18").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A19").NumberFormat = "@"
$ws.Range("A19").Value = "Critério:"
$ws.Range("A3").Copy()
$ws.Range("A19").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A20").NumberFormat = "@"
$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("A3").Copy()
$ws.Range("A20").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A21").NumberFormat = "@"
$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("A3").Copy()
$ws.Range("A21").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A22").NumberFormat = "@"
$ws.Range("A22").Value = "Requisitos:"
$ws.Range("A3").Copy()
$ws.Range("A22").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Remove cells that must no longer exist ---
$ws.Range("B14").Clear()
$ws.Range("C14").Clear()
$ws.Range("B16").Clear()
$ws.Range("C16").Clear()
$ws.Range("B22").Clear()
$ws.Range("C22").Clear()
$ws.Range("A23").Clear()

# --- Adjust row heights that differ from the target layout ---
$ws.Rows.Item(13).RowHeight = 60
$ws.Rows.Item(15).RowHeight = 120
$ws.Rows.Item(17).AutoFit()
$ws.Rows.Item(18).RowHeight = 60
$ws.Rows.Item(21).RowHeight = 120
$ws.Rows.Item(22).AutoFit()
$ws.Rows.Item(23).RowHeight = 30

# --- Delete the now-removed trailing row (former row 24) ---
$ws.Rows.Item(24).Delete()

